$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column headers I0 (I1) and IF (J1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy formatting (style) from the existing header cell H1 so the new
# headers match the bold/bordered/centered header style used across row 1.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row data: row number, I value, J value
$data = @(
    @(2, 7, 7),
    @(3, 9, 9),
    @(4, 7, 7),
    @(5, 6, 6),
    @(6, 8, 8),
    @(7, 8, 8),
    @(8, 8, 8),
    @(9, 7, 7),
    @(10, 8, 8),
    @(11, 5, 6),
    @(12, 8, 9),
    @(13, 8, 8),
    @(14, 8, 9),
    @(15, 9, 9),
    @(16, 9, 9),
    @(17, 8, 8),
    @(18, 8, 8),
    @(19, 8, 8),
    @(20, 8, 8),
    @(21, 8, 8),
    @(22, 7, 8),
    @(23, 7, 7),
    @(24, 8, 8),
    @(25, 8, 8),
    @(26, 7, 7),
    @(27, 8, 9),
    @(28, 8, 8),
    @(29, 7, 7),
    @(30, 8, 8),
    @(31, 7, 7),
    @(32, 8, 8),
    @(33, 8, 8),
    @(34, 8, 8),
    @(35, 7, 7),
    @(36, 9, 9),
    @(37, 7, 7),
    @(38, 9, 9),
    @(39, 9, 9),
    @(40, 7, 7),
    @(41, 8, 9),
    @(42, 8, 8),
    @(43, 8, 8),
    @(44, 7, 7),
    @(45, 8, 8),
    @(46, 8, 8),
    @(47, 8, 8),
    @(48, 9, 9),
    @(49, 7, 7),
    @(50, 7, 7),
    @(51, 7, 8),
    @(52, 8, 8),
    @(53, 6, 7),
    @(54, 6, 6),
    @(55, 7, 7),
    @(56, 7, 7),
    @(57, 6, 6),
    @(58, 6, 6),
    @(59, 6, 6),
    @(60, 9, 9),
    @(61, 6, 6),
    @(62, 7, 7),
    @(63, 2, 2),
    @(64, 8, 8),
    @(65, 6, 6),
    @(66, 6, 6),
    @(67, 8, 8),
    @(68, 4, 4)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
